# [FLOORPLAN_DONE] working on routing
# Update the ESP32-C6 pin mapping ("can board function") column E on Sheet1.
# The pin name/number/type/function (columns A-D) stay the same; only the
# GPIO_N routing labels in column E change (plus two renames and a TX/RX swap).
# Columns H and L mirror column E through formulas (=E.. ) and will recalc
# automatically once the source cells are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value  = "GPIO_5"
$ws.Range("E6").Value  = "GPIO_4"
$ws.Range("E7").Value  = "GPIO_3"
$ws.Range("E8").Value  = "GPIO_2"
$ws.Range("E9").Value  = "GPIO_1"
$ws.Range("E10").Value = "GPIO_0"

$ws.Range("E14").Value = "USB_N"
$ws.Range("E15").Value = "USB_P"

$ws.Range("E16").Value = "GPIO_8"
$ws.Range("E17").Value = "GPIO_8"
$ws.Range("E18").Value = "GPIO_9"
$ws.Range("E19").Value = "GPIO_10"
$ws.Range("E20").Value = "GPIO_11"
$ws.Range("E21").Value = "GPIO_12"
$ws.Range("E22").Value = "GPIO_13"
$ws.Range("E24").Value = "GPIO_14"
$ws.Range("E25").Value = "GPIO_15"
$ws.Range("E26").Value = "GPIO_6"

$ws.Range("E27").Value = "RX"
$ws.Range("E28").Value = "TX"

# Move the sheet's active selection, matching the saved view state in the
# edited workbook.
$ws.Range("E15").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
